$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.654.47"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "3.524.23"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.88%  "
$ws.Range("D7").Value = "3.524.91"
$ws.Range("E7").Value = "  -3.35%  "
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.509"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  -5.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.405"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.07%  "
$ws.Range("D13").Value = "4.120.00"
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("E14").Value = "  -7.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.80%  "
$ws.Range("D16").Value = "3.524.77"
$ws.Range("E16").Value = "  -3.06%  "
$ws.Range("D17").Value = "66.486.45"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "423.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.589"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.50%  "
$ws.Range("D25").Value = "3.673.85"
$ws.Range("E25").Value = "  -2.92%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -6.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "3.533.94"
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.153"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.46%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "173.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0810"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.60%  "
$ws.Range("E43").Value = "  -5.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  -9.45%  "
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("E49").Value = "  -5.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.904"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.14%  "
